# Update countries & provincias Spain
# - Refresh several countries' COVID stats (Alemania, Dinamarca, Trinidad y Tobago)
# - Insert/refresh "Republica Dominicana" figures, which pushes "Panama" and
#   "Luxemburgo" down a row in the table (rows 47-49)
# - Bump the "Datos actualizados" timestamp in the title cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Alemania -------------------------------------------------------
$ws.Range("B8").Value = 132553
$ws.Range("C8").Value = 343
$ws.Range("E8").Value = 56425
$ws.Range("G8").Value = 33
$ws.Range("H8").Value = 3528

# --- Row 33: Dinamarca ------------------------------------------------------
$ws.Range("D33").Value = 2748
$ws.Range("E33").Value = 3624
$ws.Range("F33").Value = 89

# --- Rows 47-49: Republica Dominicana / Panama / Luxemburgo ----------------
# Republica Dominicana now leads with freshly updated figures, and the two
# countries that used to occupy rows 47-48 (Panama, Luxemburgo) are shifted
# down one row, keeping their previous data.
$ws.Range("A47").Value = "Republica Dominicana"
$ws.Range("B47").Value = 3614
$ws.Range("C47").Value = 328
$ws.Range("D47").Value = 208
$ws.Range("E47").Value = 3217
$ws.Range("F47").Value = 143
$ws.Range("G47").Value = 6
$ws.Range("H47").Value = 189

$ws.Range("A48").Value = "Panama"
$ws.Range("B48").Value = 3574
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 72
$ws.Range("E48").Value = 3407
$ws.Range("F48").Value = 106
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 95

$ws.Range("A49").Value = "Luxemburgo"
$ws.Range("B49").Value = 3307
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 500
$ws.Range("E49").Value = 2740
$ws.Range("F49").Value = 30
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 67

# --- Row 132: Trinidad y Tobago ---------------------------------------------
$ws.Range("B132").Value = 114
$ws.Range("C132").Value = 1
$ws.Range("D132").Value = 19
$ws.Range("E132").Value = 87

# --- Title cell: refresh "last updated" timestamp --------------------------
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 16:52"
